$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1316.8889
$ws.Range("I28").Value = 731.5
$ws.Range("J28").Value = 6000
$ws.Range("K28").Value = 731.5
$ws.Range("L28").Value = 6000
$ws.Range("M28").Value = -246.5
$ws.Range("N28").Value = -6970
$ws.Range("H43").Value = 3604.4
$ws.Range("I43").Value = 2212.8572
$ws.Range("J43").Value = 6851.3335
$ws.Range("K43").Value = 2212.8572
$ws.Range("L43").Value = 6851.3335
$ws.Range("M43").Value = -2143.8572
$ws.Range("N43").Value = -6989.3335
$ws.Range("H62").Value = 3299.6667
$ws.Range("I62").Value = 3299.6667
$ws.Range("K62").Value = 3299.6667
$ws.Range("M62").Value = -2675.6667
$ws.Range("H65").Value = 3299.6667
$ws.Range("I65").Value = 3299.6667
$ws.Range("K65").Value = 16498.3335
$ws.Range("M65").Value = -13378.3335
$ws.Range("H76").Value = 7365.6665
$ws.Range("I76").Value = 3897
$ws.Range("J76").Value = 9100
$ws.Range("K76").Value = 3897
$ws.Range("L76").Value = 9100
$ws.Range("M76").Value = -3582
$ws.Range("N76").Value = -9730
$ws.Range("H79").Value = 7365.6665
$ws.Range("I79").Value = 3897
$ws.Range("J79").Value = 9100
$ws.Range("K79").Value = 3897
$ws.Range("L79").Value = 9100
$ws.Range("M79").Value = -2805
$ws.Range("N79").Value = -11284
$ws.Range("H107").Value = 1072.0555
$ws.Range("I107").Value = 934.4666999999999
$ws.Range("K107").Value = 934.4666999999999
$ws.Range("M107").Value = 985.5333000000001
$ws.Range("H113").Value = 6047.154
$ws.Range("I113").Value = 4973.7144
$ws.Range("J113").Value = 7299.5
$ws.Range("K113").Value = 4973.7144
$ws.Range("L113").Value = 7299.5
$ws.Range("M113").Value = -1719.7144
$ws.Range("N113").Value = -13807.5
$ws.Range("H129").Value = 4768
$ws.Range("I129").Value = 3456
$ws.Range("J129").Value = 4986.6665
$ws.Range("K129").Value = 10368
$ws.Range("L129").Value = 14959.9995
$ws.Range("M129").Value = -5368
$ws.Range("N129").Value = -24959.9995
$ws.Range("H141").Value = 5445.7
$ws.Range("I141").Value = 3938.1333
$ws.Range("K141").Value = 11814.3999
$ws.Range("M141").Value = -6634.3999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 557.3
$ws.Range("I4").Value = 572.125
$ws.Range("K4").Value = 572.125
$ws.Range("M4").Value = -456.125
$ws.Range("H45").Value = 2874.25
$ws.Range("I45").Value = 2550
$ws.Range("K45").Value = 2550
$ws.Range("M45").Value = -2173
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H74").Value = 1843.9231
$ws.Range("I74").Value = 1543
$ws.Range("K74").Value = 1543
$ws.Range("M74").Value = -669
$ws.Range("H77").Value = 1843.9231
$ws.Range("I77").Value = 1543
$ws.Range("K77").Value = 7715
$ws.Range("M77").Value = -3347
$ws.Range("H88").Value = 2204.0833
$ws.Range("I88").Value = 1591.6666
$ws.Range("J88").Value = 2408.2222
$ws.Range("K88").Value = 1591.6666
$ws.Range("L88").Value = 2408.2222
$ws.Range("M88").Value = -1185.6666
$ws.Range("N88").Value = -3220.2222
$ws.Range("H91").Value = 2204.0833
$ws.Range("I91").Value = 1591.6666
$ws.Range("J91").Value = 2408.2222
$ws.Range("K91").Value = 1591.6666
$ws.Range("L91").Value = 2408.2222
$ws.Range("M91").Value = -187.6666
$ws.Range("N91").Value = -5216.2222
$ws.Range("H110").Value = 5840.125
$ws.Range("I110").Value = 6245.857
$ws.Range("K110").Value = 6245.857
$ws.Range("M110").Value = -4200.857
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 381.25
$ws.Range("I94").Value = 381.25
$ws.Range("K94").Value = 381.25
$ws.Range("M94").Value = 69.75
$ws.Range("H122").Value = 38000
$ws.Range("J122").Value = 38000
$ws.Range("L122").Value = 38000
$ws.Range("N122").Value = -47800

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 5545.6294
$ws.Range("J22").Value = 9383.571
$ws.Range("L22").Value = 9383.571
$ws.Range("N22").Value = -10083.571
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").ClearContents()
$ws.Range("H132").Value = 2355.3333
$ws.Range("J132").Value = 3273.5
$ws.Range("L132").Value = 9820.5
$ws.Range("N132").Value = -14880.5
$ws.Range("H134").Value = 3903.1
$ws.Range("I134").Value = 2998.8572
$ws.Range("J134").Value = 6013
$ws.Range("K134").Value = 8996.571599999999
$ws.Range("L134").Value = 18039
$ws.Range("M134").Value = -6461.571599999999
$ws.Range("N134").Value = -23109

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5033
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 5033
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 15099
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -15323
$ws.Range("H97").Value = 343.66666
$ws.Range("J97").Value = 221
$ws.Range("L97").Value = 663
$ws.Range("N97").Value = -1655
$ws.Range("H112").Value = 9540.333000000001
$ws.Range("I112").Value = 970.75
$ws.Range("J112").Value = 11254.25
$ws.Range("K112").Value = 2912.25
$ws.Range("L112").Value = 33762.75
$ws.Range("M112").Value = -1804.25
$ws.Range("N112").Value = -35978.75
$ws.Range("H118").Value = 2998.5
$ws.Range("I118").Value = 2998.5
$ws.Range("K118").Value = 8995.5
$ws.Range("M118").Value = -7752.5
$ws.Range("H133").Value = 17210.154
$ws.Range("I133").Value = 12049.75
$ws.Range("J133").Value = 19503.666
$ws.Range("K133").Value = 36149.25
$ws.Range("L133").Value = 58510.99800000001
$ws.Range("M133").Value = -31089.25
$ws.Range("N133").Value = -68630.99800000001
$ws.Range("H137").Value = 4241.364
$ws.Range("J137").Value = 4837.2856
$ws.Range("L137").Value = 14511.8568
$ws.Range("N137").Value = -24711.8568
$ws.Range("H138").Value = 10676.667
$ws.Range("I138").Value = 10676.667
$ws.Range("K138").Value = 32030.001
$ws.Range("M138").Value = -26890.001
$ws.Range("H139").Value = 3244.4443
$ws.Range("I139").Value = 3440
$ws.Range("K139").Value = 10320
$ws.Range("M139").Value = -5180

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1950
$ws.Range("I80").Value = 1975
$ws.Range("J80").Value = 1900
$ws.Range("K80").Value = 1975
$ws.Range("L80").Value = 1900
$ws.Range("M80").Value = -977
$ws.Range("N80").Value = -3896
$ws.Range("H83").Value = 1950
$ws.Range("I83").Value = 1975
$ws.Range("J83").Value = 1900
$ws.Range("K83").Value = 9875
$ws.Range("L83").Value = 9500
$ws.Range("M83").Value = -4883
$ws.Range("N83").Value = -19484
$ws.Range("H97").Value = 765.625
$ws.Range("J97").Value = 957
$ws.Range("L97").Value = 957
$ws.Range("N97").Value = -1949
$ws.Range("H122").Value = 2172.524
$ws.Range("I122").Value = 1513.2354
$ws.Range("K122").Value = 4539.706200000001
$ws.Range("M122").Value = -2089.706200000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 9004
$ws.Range("I122").Value = 9004
$ws.Range("M122").Value = -24562
